# Auto-generated update of market-price derived columns (H-N) across the
# per-job Leve profit tables. Values come from a scheduled market-data
# refresh; only numeric result columns change, never the leve metadata.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1030.4688
$ws.Range("J17").Value = 1030.4688
$ws.Range("L17").Value = 3091.4064
$ws.Range("N17").Value = -3427.4064
$ws.Range("H107").Value = 1271.2174
$ws.Range("I107").Value = 1354.8889
$ws.Range("J107").Value = 970
$ws.Range("K107").Value = 1354.8889
$ws.Range("L107").Value = 970
$ws.Range("M107").Value = 565.1111000000001
$ws.Range("N107").Value = -4810
$ws.Range("H121").Value = 1624.674
$ws.Range("J121").Value = 1658.6666
$ws.Range("L121").Value = 4975.9998
$ws.Range("N121").Value = -8469.9998
$ws.Range("H129").Value = 1046.0817
$ws.Range("I129").Value = 449
$ws.Range("J129").Value = 1129.3954
$ws.Range("K129").Value = 1347
$ws.Range("L129").Value = 3388.1862
$ws.Range("M129").Value = 3653
$ws.Range("N129").Value = -13388.1862
$ws.Range("H132").Value = 5575.5
$ws.Range("I132").Value = 10498.889
$ws.Range("J132").Value = 652.1111
$ws.Range("K132").Value = 31496.667
$ws.Range("L132").Value = 1956.3333
$ws.Range("M132").Value = -28966.667
$ws.Range("N132").Value = -7016.3333
$ws.Range("H138").Value = 4710.778
$ws.Range("I138").Value = 3236.7307
$ws.Range("J138").Value = 6079.5356
$ws.Range("K138").Value = 9710.1921
$ws.Range("L138").Value = 18238.6068
$ws.Range("M138").Value = -4570.1921
$ws.Range("N138").Value = -28518.6068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10574.1
$ws.Range("I28").Value = 6193.3335
$ws.Range("K28").Value = 6193.3335
$ws.Range("M28").Value = -6001.3335
$ws.Range("H32").Value = 24077.684
$ws.Range("I32").Value = 18551.191
$ws.Range("J32").Value = 59999.875
$ws.Range("K32").Value = 18551.191
$ws.Range("L32").Value = 59999.875
$ws.Range("M32").Value = -18264.191
$ws.Range("N32").Value = -60573.875
$ws.Range("H37").Value = 10587.714
$ws.Range("J37").Value = 10587.714
$ws.Range("L37").Value = 10587.714
$ws.Range("N37").Value = -11133.714
$ws.Range("H99").Value = 10574.1
$ws.Range("I99").Value = 6193.3335
$ws.Range("K99").Value = 6193.3335
$ws.Range("M99").Value = -3198.3335
$ws.Range("H122").Value = 2051.3044
$ws.Range("I122").Value = 1962.8334
$ws.Range("J122").Value = 2082.5293
$ws.Range("K122").Value = 5888.5002
$ws.Range("L122").Value = 6247.5879
$ws.Range("M122").Value = -3438.5002
$ws.Range("N122").Value = -11147.5879
$ws.Range("H132").Value = 934139.8
$ws.Range("I132").Value = 1233663.2
$ws.Range("J132").Value = 2289.111
$ws.Range("K132").Value = 3700989.6
$ws.Range("L132").Value = 6867.333
$ws.Range("M132").Value = -3698459.6
$ws.Range("N132").Value = -11927.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 59780
$ws.Range("J51").Value = 59780
$ws.Range("L51").Value = 59780
$ws.Range("N51").Value = -60762

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 3333.3333
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 5000
$ws.Range("N29").Value = -5586
$ws.Range("H50").Value = 11730.667
$ws.Range("J50").Value = 11730.667
$ws.Range("L50").Value = 11730.667
$ws.Range("N50").Value = -12980.667
$ws.Range("H51").Value = 11071.429
$ws.Range("J51").Value = 11980
$ws.Range("L51").Value = 11980
$ws.Range("N51").Value = -13452
$ws.Range("H59").Value = 14398.1
$ws.Range("J59").Value = 15372.625
$ws.Range("L59").Value = 15372.625
$ws.Range("N59").Value = -17662.625
$ws.Range("H60").Value = 8873.272000000001
$ws.Range("I60").Value = 5500
$ws.Range("J60").Value = 10138.25
$ws.Range("K60").Value = 5500
$ws.Range("L60").Value = 10138.25
$ws.Range("M60").Value = -4989
$ws.Range("N60").Value = -11160.25
$ws.Range("H61").Value = 11071.429
$ws.Range("J61").Value = 11980
$ws.Range("L61").Value = 11980
$ws.Range("N61").Value = -12676
$ws.Range("H68").Value = 18385.363
$ws.Range("J68").Value = 18797.1
$ws.Range("L68").Value = 18797.1
$ws.Range("N68").Value = -20295.1
$ws.Range("H71").Value = 18385.363
$ws.Range("J71").Value = 18797.1
$ws.Range("L71").Value = 56391.3
$ws.Range("N71").Value = -63879.3
$ws.Range("H74").Value = 15207
$ws.Range("I74").Value = 9335
$ws.Range("J74").Value = 16675
$ws.Range("K74").Value = 9335
$ws.Range("L74").Value = 16675
$ws.Range("M74").Value = -8461
$ws.Range("N74").Value = -18423
$ws.Range("H77").Value = 15207
$ws.Range("I77").Value = 9335
$ws.Range("J77").Value = 16675
$ws.Range("K77").Value = 28005
$ws.Range("L77").Value = 50025
$ws.Range("M77").Value = -23637
$ws.Range("N77").Value = -58761
$ws.Range("H94").Value = 1680.5
$ws.Range("J94").Value = 1876
$ws.Range("L94").Value = 1876
$ws.Range("N94").Value = -2778
$ws.Range("H99").Value = 112346
$ws.Range("I99").Value = 1328.5714
$ws.Range("J99").Value = 500907
$ws.Range("K99").Value = 1328.5714
$ws.Range("L99").Value = 500907
$ws.Range("M99").Value = 169.4286
$ws.Range("N99").Value = -503903
$ws.Range("H126").Value = 112346
$ws.Range("I126").Value = 1328.5714
$ws.Range("J126").Value = 500907
$ws.Range("K126").Value = 3985.7142
$ws.Range("L126").Value = 1502721
$ws.Range("M126").Value = -1515.7142
$ws.Range("N126").Value = -1507661
$ws.Range("H132").Value = 1631.0278
$ws.Range("I132").Value = 1327.875
$ws.Range("J132").Value = 4056.25
$ws.Range("K132").Value = 3983.625
$ws.Range("L132").Value = 12168.75
$ws.Range("M132").Value = -1453.625
$ws.Range("N132").Value = -17228.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 1536.25
$ws.Range("J94").Value = 2510
$ws.Range("L94").Value = 7530
$ws.Range("N94").Value = -8882
$ws.Range("H106").Value = 5559889
$ws.Range("J106").Value = 5559889
$ws.Range("L106").Value = 16679667
$ws.Range("N106").Value = -16681559
$ws.Range("H113").Value = 636.9167
$ws.Range("I113").Value = 783
$ws.Range("J113").Value = 588.2222
$ws.Range("K113").Value = 2349
$ws.Range("L113").Value = 1764.6666
$ws.Range("M113").Value = -179
$ws.Range("N113").Value = -6104.6666
$ws.Range("H131").Value = 44876656
$ws.Range("I131").Value = 133345430
$ws.Range("J131").Value = 23812662
$ws.Range("K131").Value = 400036290
$ws.Range("L131").Value = 71437986
$ws.Range("M131").Value = -400031250
$ws.Range("N131").Value = -71448066
$ws.Range("H140").Value = 2020.3448
$ws.Range("I140").Value = 1510
$ws.Range("J140").Value = 2990
$ws.Range("K140").Value = 4530
$ws.Range("L140").Value = 8970
$ws.Range("M140").Value = 650
$ws.Range("N140").Value = -19330

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2733.5715
$ws.Range("I97").Value = 2756.6667
$ws.Range("J97").Value = 2692
$ws.Range("K97").Value = 2756.6667
$ws.Range("L97").Value = 2692
$ws.Range("M97").Value = -2260.6667
$ws.Range("N97").Value = -3684
$ws.Range("H122").Value = 6367.25
$ws.Range("I122").Value = 8586.714
$ws.Range("K122").Value = 25760.142
$ws.Range("M122").Value = -23310.142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1659.5333
$ws.Range("I16").Value = 1792.3704
$ws.Range("J16").Value = 464
$ws.Range("K16").Value = 1792.3704
$ws.Range("L16").Value = 464
$ws.Range("M16").Value = -1622.3704
$ws.Range("N16").Value = -804
$ws.Range("H40").Value = 1827.5
$ws.Range("I40").Value = 2350
$ws.Range("J40").Value = 1305
$ws.Range("K40").Value = 2350
$ws.Range("L40").Value = 1305
$ws.Range("M40").Value = -2214
$ws.Range("N40").Value = -1577
$ws.Range("H46").Value = 1013.5455
$ws.Range("I46").Value = 1015.7895
$ws.Range("J46").Value = 999.3333
$ws.Range("K46").Value = 1015.7895
$ws.Range("L46").Value = 999.3333
$ws.Range("M46").Value = -827.7895
$ws.Range("N46").Value = -1375.3333
$ws.Range("H122").Value = 2060.9387
$ws.Range("I122").Value = 1976.1316
$ws.Range("J122").Value = 2353.9092
$ws.Range("K122").Value = 5928.3948
$ws.Range("L122").Value = 7061.7276
$ws.Range("M122").Value = -3478.3948
$ws.Range("N122").Value = -11961.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 55010
$ws.Range("J24").Value = 55010
$ws.Range("L24").Value = 55010
$ws.Range("N24").Value = -55470
$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40452
